$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Edit 1: "Постановка задачи: ... подпроцессов." paragraph gains
#         " Вариант 30." appended after the final period.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(52)
$pr1 = $p1.Range
$start1 = $pr1.End - 2
$end1 = $pr1.End - 1
$r1 = $d.Range($start1, $end1)
$xml1 = '<w:p ' + $wns + '>' `
  + '<w:r w:rsidR="002E2E5B" w:rsidRPr="002E2E5B"><w:rPr><w:lang w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>.</w:t></w:r>' `
  + '<w:r><w:rPr><w:lang w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t xml:space="preserve"> Вариант 30</w:t></w:r>' `
  + '<w:r><w:rPr><w:lang w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>.</w:t></w:r>' `
  + '</w:p>'
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: Figure 1 caption "Производство икры A0" -> "Контекстная диаграмма"
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(55)
$pr2 = $p2.Range
$t2 = $pr2.Text
$idx2 = $t2.IndexOf("Производство")
$start2 = $pr2.Start + $idx2
$end2 = $pr2.End - 1
$r2 = $d.Range($start2, $end2)
$xml2 = '<w:p ' + $wns + '><w:r w:rsidR="0096556E"><w:rPr><w:lang w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>Контекстная диаграмма</w:t></w:r></w:p>'
$r2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Edit 3: Figure 2 caption "Производство икры A1, А2, А3" ->
#         "Детализация контекстной диаграммы"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(57)
$pr3 = $p3.Range
$t3 = $pr3.Text
$idx3 = $t3.IndexOf("Производство")
$start3 = $pr3.Start + $idx3
$end3 = $pr3.End - 1
$r3 = $d.Range($start3, $end3)
$xml3 = '<w:p ' + $wns + '><w:r w:rsidR="0096556E"><w:rPr><w:lang w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>Детализация контекстной диаграммы</w:t></w:r></w:p>'
$r3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Edit 4: Figure 3 caption "Производство икры A21, А22, А23" ->
#         "Детализация одного из процессов"
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(59)
$pr4 = $p4.Range
$t4 = $pr4.Text
$idx4 = $t4.IndexOf("Производство")
$start4 = $pr4.Start + $idx4
$end4 = $pr4.End - 1
$r4 = $d.Range($start4, $end4)
$xml4 = '<w:p ' + $wns + '><w:r w:rsidR="0096556E"><w:rPr><w:lang w:eastAsia="zh-CN" w:bidi="hi-IN"/></w:rPr><w:t>Детализация одного из процессов</w:t></w:r></w:p>'
$r4.InsertXML($xml4)

Write-Host "Done"
